# Apply the "address label regex" changes to Sheet1 of the LiveConfig workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 53 currently holds "SeparateProducts" in A53:C53; the table ends at row 54
# which is a mostly-empty row (just A54 with a style, no value).
# We need to turn that trailing empty row into real data, and append one more row
# after it, growing Table1 from A1:C54 to A1:C55.

# Copy the formatting of the preceding data row (53) down onto the two new
# rows so the new cells pick up the normal (non-header) "Name"/"Value" look
# (style indices 3/4) rather than leaving row 54's old placeholder header
# style (index 5) in place.
$ws.Range("A53:B53").Copy()
$ws.Range("A54:B55").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(54, 1).Value = "AddressContinueRegex"
$ws.Cells.Item(54, 2).Value = "\WContinue\W:(\w+)"

$ws.Cells.Item(55, 1).Value = "AddressRegex"
$ws.Cells.Item(55, 2).Value = 'IN_strAddressLabel\W{0,2}:\W{0,2}(.+)",'

# Grow the table (ListObject) so it covers the newly-added row 55 as well.
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:C55"))

$ws.Range("B57").Select()
